$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a brand-new row at position 13 (pushes the old rows 13-21 down to
#    14-22). The new row only has content in columns B and C (the docente
#    responsible), column A stays completely empty (no cell at all), and the
#    row keeps the default height.
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Clear()

# Grab the formatting (font/alignment) that column B/C normally use - the
# shifted-down row 14 still carries it at this point - and apply it to the
# new row 13 before we overwrite row 14's text.
$ws.Range("B14:C14").Copy() | Out-Null
$ws.Range("B13:C13").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("B13").Value = "198273 - Domingos Savio Giordani"
$ws.Range("C13").Value = "198273 - Domingos Savio Giordani"

# ---------------------------------------------------------------------------
# 2. Fix up the content that was wrong / incomplete in the original sheet.
# ---------------------------------------------------------------------------

# Row 10 - Objetivos / Objectives: B10/C10 mistakenly held the "Docentes" text
$objetivos = "Este curso tem por objetivo fornecer aos alunos de Engenharia de Produção os princípios fundamentais da Química com enfoque tecnológico e nas aplicações industriais passíveis de serem encontradas na profissão."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# Row 14 - Programa resumido: previously just "Semestral"
$programaResumido = "1 – Conceitos básicos de Química; 2 – Os estados físicos da matéria e suas propriedades peculiares; 3 – Reações químicas; 4 – Noções de química orgânica; 5 – Materiais modernos; 6 – Tecnologia Química aplicada"
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido

# Row 16 - Programa: previously mistakenly held a date
$programa = "Programa em português1.Conceitos básicos de Química (2 horas)a.Estrutura Atômicab.Tabela Periódicac.Ligações Químicas2.Os estados físicos da matéria e suas propriedades peculiares (6 horas)a.O estado gasoso – pressão, relações PVT, gases ideais e reaisb.O estado líquido – soluções, forças intermoleculares, viscosidade, tensão superficial, pressão de vapor, mudanças de fasec.O estado sólido – classificação dos sólidos (moleculares, reticulares, metálicos e iônicos) 3.Reações químicas (8 horas)a.Tipos de reações (dupla-troca, oxirredução)b.Estequiometria em reações químicas (reagentes limitantes, pureza e rendimento)c.Energia e reações químicasd.Equilíbrio químico – soluções tampãoe.Fundamentos de corrosão4.Noções de química orgânica (6 horas)a.Hidrocarbonetos e suas principais propriedadesb.Combustíveis e combustãoc.Polímeros5.Tecnologia Química aplicada (8 horas)a.Papel e celuloseb.Açúcar e álcoolc.Sabões e detergentesd.Petróleo e gáse.Gases industriais f.Produção de vidros e cimento"
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# Row 19 - Método: previously mistakenly held the "Docentes" text
$metodo = "Duas provas escritas e um seminário que, juntos, constituem a primeira avaliação."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# Row 20 - Critério: gets the grading-average text that used to sit one row up
$criterio = "A nota de primeira avaliação será igual à média das notas das duas provas, com peso 7 somada à nota do seminário com peso 3. Alunos com nota de primeira avaliação igual ou superior a 5 estarão aprovados, com nota entre 3 e 4,9 em recuperação e abaixo de 3 reprovados."
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# Row 21 - Norma de recuperação: gets the recovery-exam text that used to sit one row up
$normaRecuperacao = "A recuperação se constituirá de uma prova abordando todos os assuntos do semestre, a nota de segunda avaliação será igual à média entre a nota de primeira avaliação e a prova de recuperação. Alunos com nota de segunda avaliação igual ou superior a 5 estarão aprovados e inferior a 5 reprovados."
$ws.Range("B21").Value = $normaRecuperacao
$ws.Range("C21").Value = $normaRecuperacao

# Row 22 - Bibliografia: previously empty, now filled in
$bibliografia = "BROWN, T.L. et al. Química a ciência central. 9.ed. São Paulo: Pearson Prentice Hall, 2005-2007ATKINS, P. Princípios de Química, questionando a vida moderna e o meio ambiente. 3ª Ed. Porto Alegre: Editora Bookman, 2006KOTZ, J. C. et al. Química geral e reações químicas, 9ª Edição, São Paulo, Cengage Learning, 2015.TOLENTINO, N. M. C. Processos Químicos Industriais, 1ª Edição, São Paulo, Érica, 2015."
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia

# ---------------------------------------------------------------------------
# 3. Column A no longer shares its width/style definition with column B;
#    touching column B's width makes the workbook store the two column
#    ranges separately again, while column A keeps its original width.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(2).ColumnWidth
